$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert a new row at position 389 (shifts existing rows 389..448 down to 390..449)
$ws.Rows.Item(389).Insert()

# Populate the newly inserted row 389 with its data
$ws.Cells.Item(389, 1).Value = 9
$ws.Cells.Item(389, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(389, 3).Value = "Metropolitana"
$ws.Cells.Item(389, 4).Value = 45180
$ws.Cells.Item(389, 5).Value = 13
$ws.Cells.Item(389, 6).Value = 100112021
$ws.Cells.Item(389, 7).Value = "Ají"
$ws.Cells.Item(389, 8).Value = "Inferno"
$ws.Cells.Item(389, 9).Value = "Primera"
$ws.Cells.Item(389, 10).Value = 52
$ws.Cells.Item(389, 11).Value = 29000
$ws.Cells.Item(389, 12).Value = 30000
$ws.Cells.Item(389, 13).Value = 29500
$ws.Cells.Item(389, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(389, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(389, 16).Value = 2950
$ws.Cells.Item(389, 17).Value = 10
$ws.Cells.Item(389, 18).Value = "Hortaliza"
